$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.936.71'
$ws.Range("E2").Value = '  +1.95%  '

$ws.Range("D3").Value = '3.110.01'
$ws.Range("E3").Value = '  +5.37%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.08'
$ws.Range("E5").Value = '  +1.58%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.08'
$ws.Range("E6").Value = '  +6.57%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").Value = '3.106.02'
$ws.Range("E8").Value = '  +5.38%  '

$ws.Range("E9").Value = '  +1.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.49'
$ws.Range("E10").Value = '  -3.59%  '

$ws.Range("E11").Value = '  +3.53%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.485'
$ws.Range("E12").Value = '  +5.17%  '

$ws.Range("E13").Value = '  +1.70%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.53'
$ws.Range("E14").Value = '  +7.80%  '

$ws.Range("E15").Value = '  -0.06%  '

$ws.Range("D16").Value = '3.626.33'
$ws.Range("E16").Value = '  +5.32%  '

$ws.Range("D17").Value = '66.902.95'
$ws.Range("E17").Value = '  +1.91%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.20'
$ws.Range("E18").Value = '  +2.61%  '

$ws.Range("D19").Value = '3.105.23'
$ws.Range("E19").Value = '  +5.19%  '

$ws.Range("E20").Value = '  +1.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '484.85'
$ws.Range("E21").Value = '  +8.79%  '

$ws.Range("E22").Value = '  +2.59%  '

$ws.Range("E23").Value = '  +3.26%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.18'
$ws.Range("E24").Value = '  +2.42%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.40'
$ws.Range("E25").Value = '  +6.38%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.16'
$ws.Range("E26").Value = '  +7.18%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.08'
$ws.Range("E27").Value = '  +0.48%  '

$ws.Range("E28").Value = '  -0.08%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.00'
$ws.Range("E29").Value = '  -1.19%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.41'
$ws.Range("E30").Value = '  -4.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.70'
$ws.Range("E31").Value = '  +3.92%  '

$ws.Range("B32").Value = 'PEPE'
$ws.Range("C32").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0000101'
$ws.Range("E32").Value = '  -1.82%  '

$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.00'
$ws.Range("E33").Value = '  +6.35%  '

$ws.Range("E34").Value = '  +1.63%  '

$ws.Range("E35").Value = '  -0.10%  '

$ws.Range("B36").Value = 'Mantle'
$ws.Range("C36").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +3.40%  '

$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.92'
$ws.Range("E37").Value = '  +3.38%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '47.72'
$ws.Range("E38").Value = '  +4.55%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.14'
$ws.Range("E39").Value = '  +7.39%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.25'
$ws.Range("E40").Value = '  +2.18%  '

$ws.Range("E41").Value = '  +4.87%  '

$ws.Range("E42").Value = '  +0.51%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.69'
$ws.Range("E43").Value = '  +1.54%  '

$ws.Range("E44").Value = '  -1.37%  '

$ws.Range("E45").Value = '  +2.68%  '

$ws.Range("D46").Value = '2.842.02'
$ws.Range("E46").Value = '  +5.96%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '383.38'
$ws.Range("E47").Value = '  -0.87%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '135.02'
$ws.Range("E48").Value = '  +1.10%  '

$ws.Range("E49").Value = '  -0.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.02'
$ws.Range("E50").Value = '  +4.92%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.24'
$ws.Range("E51").Value = '  +3.07%  '
